$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from A48 to the new A49 cell before setting values,
# so the new row matches the existing styled "index" column.
$ws.Cells.Item(48, 1).Copy()
$ws.Cells.Item(49, 1).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New index value for the appended row (A49 = 47)
$ws.Cells.Item(49, 1).Value = 47

# Updated time-series values (columns B and C) for rows 2-49,
# reflecting the re-run simulation with the new priority-queue based algorithm.
$ws.Cells.Item(2, 2).Value = 3.39105785399794
$ws.Cells.Item(2, 3).Value = 0.6992512683249379
$ws.Cells.Item(3, 2).Value = 4.776038466560848
$ws.Cells.Item(3, 3).Value = 1.487991444271757
$ws.Cells.Item(4, 2).Value = 4.906956361660531
$ws.Cells.Item(4, 3).Value = 2.346843977729143
$ws.Cells.Item(5, 2).Value = 7.095214515952791
$ws.Cells.Item(5, 3).Value = 3.306489725452458
$ws.Cells.Item(6, 2).Value = 9.409135020118875
$ws.Cells.Item(6, 3).Value = 4.276037286930864
$ws.Cells.Item(7, 2).Value = 9.829943325960574
$ws.Cells.Item(7, 3).Value = 5.314247456939732
$ws.Cells.Item(8, 2).Value = 13.6996308890176
$ws.Cells.Item(8, 3).Value = 5.967354696538082
$ws.Cells.Item(9, 2).Value = 13.82085137585603
$ws.Cells.Item(9, 3).Value = 6.804923809010105
$ws.Cells.Item(10, 2).Value = 14.83491456537051
$ws.Cells.Item(10, 3).Value = 7.909190846096122
$ws.Cells.Item(11, 2).Value = 14.95911033097814
$ws.Cells.Item(11, 3).Value = 8.994304127790596
$ws.Cells.Item(12, 2).Value = 21.6965379458027
$ws.Cells.Item(12, 3).Value = 9.885675018990186
$ws.Cells.Item(13, 2).Value = 25.5421900738598
$ws.Cells.Item(13, 3).Value = 10.78153933301147
$ws.Cells.Item(14, 2).Value = 25.6218956729883
$ws.Cells.Item(14, 3).Value = 11.58945225760314
$ws.Cells.Item(15, 2).Value = 25.70805778311915
$ws.Cells.Item(15, 3).Value = 12.79137840339882
$ws.Cells.Item(16, 2).Value = 27.38371103004833
$ws.Cells.Item(16, 3).Value = 13.6459330165744
$ws.Cells.Item(17, 2).Value = 28.63990918739093
$ws.Cells.Item(17, 3).Value = 14.43633139326868
$ws.Cells.Item(18, 2).Value = 29.63641111753113
$ws.Cells.Item(18, 3).Value = 15.4657912796883
$ws.Cells.Item(19, 2).Value = 29.75427616056376
$ws.Cells.Item(19, 3).Value = 16.0315698557942
$ws.Cells.Item(20, 2).Value = 29.89312026977975
$ws.Cells.Item(20, 3).Value = 17.06009056598444
$ws.Cells.Item(21, 2).Value = 31.00338641297488
$ws.Cells.Item(21, 3).Value = 18.46590478413921
$ws.Cells.Item(22, 2).Value = 35.39749516784784
$ws.Cells.Item(22, 3).Value = 19.59962713818129
$ws.Cells.Item(23, 2).Value = 35.55091013240617
$ws.Cells.Item(23, 3).Value = 20.61625198156935
$ws.Cells.Item(24, 2).Value = 36.7207941949767
$ws.Cells.Item(24, 3).Value = 21.21338889668544
$ws.Cells.Item(25, 2).Value = 36.84376913736955
$ws.Cells.Item(25, 3).Value = 22.10420018242315
$ws.Cells.Item(26, 2).Value = 38.36621303754708
$ws.Cells.Item(26, 3).Value = 22.74893893549961
$ws.Cells.Item(27, 2).Value = 38.40951430118925
$ws.Cells.Item(27, 3).Value = 23.4104987530423
$ws.Cells.Item(28, 2).Value = 38.49666784215022
$ws.Cells.Item(28, 3).Value = 24.44732715254547
$ws.Cells.Item(29, 2).Value = 42.49994436860529
$ws.Cells.Item(29, 3).Value = 26.07089876039941
$ws.Cells.Item(30, 2).Value = 43.08412604696379
$ws.Cells.Item(30, 3).Value = 26.9133378495347
$ws.Cells.Item(31, 2).Value = 49.14814027050019
$ws.Cells.Item(31, 3).Value = 27.88392528763817
$ws.Cells.Item(32, 2).Value = 49.31310077743742
$ws.Cells.Item(32, 3).Value = 28.67426184274756
$ws.Cells.Item(33, 2).Value = 51.61157168472089
$ws.Cells.Item(33, 3).Value = 29.59611614964999
$ws.Cells.Item(34, 2).Value = 51.72711240467387
$ws.Cells.Item(34, 3).Value = 30.41291300365097
$ws.Cells.Item(35, 2).Value = 53.9077547417444
$ws.Cells.Item(35, 3).Value = 31.17816367291812
$ws.Cells.Item(36, 2).Value = 55.367479528487
$ws.Cells.Item(36, 3).Value = 31.97812058888696
$ws.Cells.Item(37, 2).Value = 58.82308468637991
$ws.Cells.Item(37, 3).Value = 32.78373362118491
$ws.Cells.Item(38, 2).Value = 58.9228285046602
$ws.Cells.Item(38, 3).Value = 33.615534256519
$ws.Cells.Item(39, 2).Value = 59.14909276317007
$ws.Cells.Item(39, 3).Value = 34.61119128681807
$ws.Cells.Item(40, 2).Value = 66.33440682034083
$ws.Cells.Item(40, 3).Value = 35.67407127430278
$ws.Cells.Item(41, 2).Value = 66.48914716532038
$ws.Cells.Item(41, 3).Value = 36.42335496790051
$ws.Cells.Item(42, 2).Value = 70.25143776101797
$ws.Cells.Item(42, 3).Value = 37.67316715366122
$ws.Cells.Item(43, 2).Value = 78.52943680029652
$ws.Cells.Item(43, 3).Value = 38.97828985616626
$ws.Cells.Item(44, 2).Value = 85.39317880388944
$ws.Cells.Item(44, 3).Value = 39.83401648375021
$ws.Cells.Item(45, 2).Value = 87.59798209938131
$ws.Cells.Item(45, 3).Value = 41.03268199820609
$ws.Cells.Item(46, 2).Value = 87.7979391329324
$ws.Cells.Item(46, 3).Value = 41.71246517789266
$ws.Cells.Item(47, 2).Value = 87.91052813281124
$ws.Cells.Item(47, 3).Value = 42.58120238901444
$ws.Cells.Item(48, 2).Value = 88.75645137962633
$ws.Cells.Item(48, 3).Value = 43.42239086856442
$ws.Cells.Item(49, 2).Value = 90.03971422891796
$ws.Cells.Item(49, 3).Value = 44.28760029239231

Write-Output "Done updating Energy Consumption data through row 49."
